# edit c_click --None issue
# Update quantity ("수량") values in the C column across several sheets.

$wb = $excel.ActiveWorkbook

# 식당판매 (sheet1): C6 0 -> 3
$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C6").Value = 3

# 매점판매 (sheet2): C4 1 -> 4, C6 5 -> 9
$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C4").Value = 4
$ws2.Range("C6").Value = 9

# 장의용품 (sheet3): add C5 = 11, clear C8 (was 0)
$ws3 = $wb.Worksheets.Item("장의용품")
$ws3.Range("C5").Value = 11
$ws3.Range("C8").ClearContents()

# 상복 (sheet4): C7 2 -> 0, C11 2 -> 9
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C7").Value = 0
$ws4.Range("C11").Value = 9

# 기타 (sheet5): C2..C14 quantities updated
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C2").Value = 3
$ws5.Range("C3").Value = 2
$ws5.Range("C4").Value = 1
$ws5.Range("C5").Value = 4
$ws5.Range("C6").Value = 5
$ws5.Range("C7").Value = 23
$ws5.Range("C8").Value = 65
$ws5.Range("C9").Value = 43
$ws5.Range("C10").Value = 233
$ws5.Range("C11").Value = 44
$ws5.Range("C12").Value = 23
$ws5.Range("C13").Value = 1
$ws5.Range("C14").Value = 43
